$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: Earnings Before Interest And Taxes - J (2011) becomes "NA"
$ws.Range("J21").Value = "NA"

# Row 83: Depreciation - J (2011) becomes "NA"
$ws.Range("J83").Value = "NA"

# Row 91: Capital Expenditures - all years (D:J) updated with new values
$ws.Range("D91").Value = -300
$ws.Range("E91").Value = -200
$ws.Range("F91").Value = -1400
$ws.Range("G91").Value = -700
$ws.Range("H91").Value = -300
$ws.Range("I91").Value = -1500
$ws.Range("J91").Value = -3100

# Row 94: Total Cash Flows From Investing Activities - J (2011) becomes "NA"
$ws.Range("J94").Value = "NA"

# Row 100: Total Cash Flows From Financing Activities - J (2011) becomes "NA"
$ws.Range("J100").Value = "NA"

# Row 101: Effect Of Exchange Rate Changes - J (2011) becomes "NA"
$ws.Range("J101").Value = "NA"
